$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.452.25'
$ws.Range("E2").Value = '  -1.60%  '
$ws.Range("D3").Value = '2.507.32'
$ws.Range("E3").Value = '  -4.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.67'
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.63'
$ws.Range("E6").Value = '  +2.31%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -2.49%  '
$ws.Range("D9").Value = '2.506.03'
$ws.Range("E9").Value = '  -4.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.349'
$ws.Range("E12").Value = '  -4.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.09'
$ws.Range("E13").Value = '  -2.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.49'
$ws.Range("E14").Value = '  -4.36%  '
$ws.Range("D15").Value = '2.939.55'
$ws.Range("E15").Value = '  -5.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000175'
$ws.Range("E16").Value = '  -4.11%  '
$ws.Range("D17").Value = '66.211.06'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").Value = '2.507.95'
$ws.Range("E18").Value = '  -5.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.23'
$ws.Range("E19").Value = '  -6.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.65'
$ws.Range("E20").Value = '  -5.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '345.90'
$ws.Range("E21").Value = '  -3.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.19'
$ws.Range("E22").Value = '  -3.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.58'
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.93'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.33'
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.83'
$ws.Range("E27").Value = '  -4.78%  '
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").Value = '2.627.21'
$ws.Range("E29").Value = '  -4.78%  '
$ws.Range("D30").Value = '0.0₃0969'
$ws.Range("E30").Value = '  -3.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '527.57'
$ws.Range("E31").Value = '  -3.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.07'
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.32'
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.83'
$ws.Range("E34").Value = '  -3.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.131'
$ws.Range("E35").Value = '  -3.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.75'
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.45'
$ws.Range("E38").Value = '  -3.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.52'
$ws.Range("E39").Value = '  -2.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.35'
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.353'
$ws.Range("E41").Value = '  -3.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.78'
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.08'
$ws.Range("E43").Value = '  -3.02%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '146.92'
$ws.Range("E46").Value = '  -3.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.555'
$ws.Range("E47").Value = '  -4.35%  '
$ws.Range("E48").Value = '  -3.48%  '
$ws.Range("E49").Value = '  +1.34%  '
$ws.Range("D50").Value = '0.0₆0270'
$ws.Range("E50").Value = '  -9.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0751'
$ws.Range("E51").Value = '  -2.53%  '
